$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 214, pushing the existing data
# (previously rows 214-297) down to rows 216-299.
$ws.Rows("214:215").Insert()

# Populate the first new row (214) with a new price record for
# "Kurakata" peaches sold at "Región de O'Higgins".
$ws.Cells.Item(214, 1).Value  = 7
$ws.Cells.Item(214, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(214, 3).Value  = "Ñuble"
$ws.Cells.Item(214, 4).Value  = 44917
$ws.Cells.Item(214, 5).Value  = 16
$ws.Cells.Item(214, 6).Value  = "Fruta"
$ws.Cells.Item(214, 7).Value  = 100103
$ws.Cells.Item(214, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(214, 9).Value  = 100103004
$ws.Cells.Item(214, 10).Value = "Durazno"
$ws.Cells.Item(214, 11).Value = "Kurakata"
$ws.Cells.Item(214, 12).Value = "Primera"
$ws.Cells.Item(214, 13).Value = 120
$ws.Cells.Item(214, 14).Value = 16000
$ws.Cells.Item(214, 15).Value = 17000
$ws.Cells.Item(214, 16).Value = 16500
$ws.Cells.Item(214, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(214, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(214, 19).Value = 1100
$ws.Cells.Item(214, 20).Value = 15

# Populate the second new row (215) with a matching record for the
# lower "Segunda" quality grade.
$ws.Cells.Item(215, 1).Value  = 7
$ws.Cells.Item(215, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(215, 3).Value  = "Ñuble"
$ws.Cells.Item(215, 4).Value  = 44917
$ws.Cells.Item(215, 5).Value  = 16
$ws.Cells.Item(215, 6).Value  = "Fruta"
$ws.Cells.Item(215, 7).Value  = 100103
$ws.Cells.Item(215, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(215, 9).Value  = 100103004
$ws.Cells.Item(215, 10).Value = "Durazno"
$ws.Cells.Item(215, 11).Value = "Kurakata"
$ws.Cells.Item(215, 12).Value = "Segunda"
$ws.Cells.Item(215, 13).Value = 60
$ws.Cells.Item(215, 14).Value = 15000
$ws.Cells.Item(215, 15).Value = 15000
$ws.Cells.Item(215, 16).Value = 15000
$ws.Cells.Item(215, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(215, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(215, 19).Value = 1000
$ws.Cells.Item(215, 20).Value = 15
